$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Release Plan")

# --- Shift Sprint 1-5 start dates / durations (Increment Plan rows 16-20) ---
# Sprint 1 (row 16): start 9/9/2021 -> 11/1/2021, duration 8 -> 7 days
$ws.Range("B16").Value = 44501
$ws.Range("C16").Value = 7
$ws.Range("H16").Value = 44507

# Sprint 2 (row 17): start 9/17/2021 -> 11/8/2021, duration 8 -> 7 days
$ws.Range("B17").Value = 44508
$ws.Range("C17").Value = 7

# Sprint 3 (row 18): start 9/25/2021 -> 11/15/2021, duration 8 -> 7 days
$ws.Range("B18").Value = 44515
$ws.Range("C18").Value = 7

# Sprint 4 (row 19): start 10/3/2021 -> 11/22/2021 (duration stays 6 days)
$ws.Range("B19").Value = 44522

# Sprint 5 (row 20): start 10/9/2021 -> 11/28/2021 (duration stays 7 days)
$ws.Range("B20").Value = 44528

# --- Release Date for Increment 1 (H4) becomes literal text "7/22/2021" ---
# instead of a date serial, while keeping the existing cell style (s="10").
# Assigning the literal string directly gets auto-coerced back into a date
# serial by the engine, so route it through a temporary formula and then
# freeze the computed text back down to a plain value via copy/paste.
$ws.Range("H4").Formula = '="7/22/2021"'
$ws.Range("H4").Copy()
$ws.Range("H4").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$excel.CutCopyMode = $false

# --- Update the active selection to match the saved view state ---
$ws.Activate()
$ws.Range("I5").Select()
